$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new day row (row 47) following the existing table pattern.
$newRow = 47

$ws.Cells.Item($newRow, 1).Value = 45996
$ws.Cells.Item($newRow, 2).Value = 107
$ws.Cells.Item($newRow, 3).Value = 119
$ws.Cells.Item($newRow, 4).Value = 114

# Match the date-formatted style used by the rest of column A (row 46).
$ws.Range("A46").Copy()
$ws.Range("A$newRow").PasteSpecial(-4122) # xlPasteFormats
